$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (21) -- shifts old U/V/W -> V/W/X
$ws.Columns("U").Insert()

# New header for the inserted column
$ws.Range("U1").Value = "csim"

# New column values for the inserted "csim" column (rows 2-11)
$csimValues = @(0, 0, 0, 1, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $csimValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 21).Value = $csimValues[$i]
}

# All "position" values (column B) become -1
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 2).Value = -1
}
